# Delete the 4th slide (an empty/blank slide with no content) from the
# presentation, matching the upstream commit "Update Sprint Cadence
# Diagrams.pptx" which removed the <p:sldId id="259" r:id="rId5"/> entry
# and its corresponding ppt/slides/slide4.xml part.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$s.Delete()
